# Sound assets and FMOD upload
# Adds a "Status" column (D) to the asset list:
#  - D1 header "Status" (styled like the other headers, with a new orange fill)
#  - D4:D18 marked "Event Created in FMOD" (note: D2/D3 intentionally stay blank,
#    matching the original author's edit)
#  - New asset row 18: White Noise / Ambience / dl from freesound / Event Created in FMOD

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Status" column header ---------------------------------------------
# (shared-string allocation order matters for an exact XML match, so this is
#  written first: Status, then White Noise / dl from freesound, then the
#  repeated "Event Created in FMOD" entries - matching the original author's
#  save order.)
$ws.Range("D1").Value = "Status"

# --- New row 18: White Noise asset -----------------------------------------
$ws.Range("A18").Value = "White Noise"
$ws.Range("B18").Value = "Ambience"
$ws.Range("C18").Value = "dl from freesound"
$ws.Range("D18").Value = "Event Created in FMOD"

$ws.Range("D4").Value = "Event Created in FMOD"
$ws.Range("D5").Value = "Event Created in FMOD"
$ws.Range("D6").Value = "Event Created in FMOD"
$ws.Range("D7").Value = "Event Created in FMOD"
$ws.Range("D8").Value = "Event Created in FMOD"
$ws.Range("D9").Value = "Event Created in FMOD"
$ws.Range("D10").Value = "Event Created in FMOD"
$ws.Range("D11").Value = "Event Created in FMOD"
$ws.Range("D12").Value = "Event Created in FMOD"
$ws.Range("D13").Value = "Event Created in FMOD"
$ws.Range("D14").Value = "Event Created in FMOD"
$ws.Range("D15").Value = "Event Created in FMOD"
$ws.Range("D16").Value = "Event Created in FMOD"
$ws.Range("D17").Value = "Event Created in FMOD"

# --- Header styling: match the other header cells (bold Heading 1 font +
#     thick bottom border) but with a new custom fill color -----------------
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D1").Interior.Color = 10079487

# --- Column width for the new column ----------------------------------------
$ws.Columns("D").ColumnWidth = 37.45

# --- Selection / view state, matching the saved workbook --------------------
$ws.Range("D2").Select()

# --- Page setup (orientation explicitly set to portrait) -------------------
$ws.PageSetup.Orientation = 1
